# Apply NATMI Cd14-Itga4 data update (Dr Hou advice revision)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 108.3097613333333
$ws.Range("H2").Value = 324.929284
$ws.Range("I2").Value = 0.9760647858278649
$ws.Range("J2").Value = 0.9760647858278649
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 26.04517333333333
$ws.Range("N2").Value = 78.13552
$ws.Range("O2").Value = 0.9210237118384171
$ws.Range("P2").Value = 0.921023711838417
$ws.Range("Q2").Value = 2820.946507618631
$ws.Range("R2").Value = 25388.51856856768
$ws.Range("S2").Value = 0.8989788120379498
$ws.Range("T2").Value = 0.8989788120379497

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 108.3097613333333
$ws.Range("H3").Value = 324.929284
$ws.Range("I3").Value = 0.9760647858278649
$ws.Range("J3").Value = 0.9760647858278649
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3302223333333333
$ws.Range("N3").Value = 0.990667
$ws.Range("O3").Value = 0.01167750336256582
$ws.Range("P3").Value = 0.01167750336256582
$ws.Range("Q3").Value = 35.76630211026978
$ws.Range("R3").Value = 321.896718992428
$ws.Range("S3").Value = 0.01139799981858698
$ws.Range("T3").Value = 0.01139799981858698

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 108.3097613333333
$ws.Range("H4").Value = 324.929284
$ws.Range("I4").Value = 0.9760647858278649
$ws.Range("J4").Value = 0.9760647858278649
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.903109
$ws.Range("N4").Value = 5.709327
$ws.Range("O4").Value = 0.06729878479901708
$ws.Range("P4").Value = 0.06729878479901708
$ws.Range("Q4").Value = 206.1252815813187
$ws.Range("R4").Value = 1855.127534231868
$ws.Range("S4").Value = 0.06568797397132817
$ws.Range("T4").Value = 0.06568797397132817

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.160250666666667
$ws.Range("H5").Value = 6.480752000000001
$ws.Range("I5").Value = 0.0194677245922947
$ws.Range("J5").Value = 0.0194677245922947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.04517333333333
$ws.Range("N5").Value = 78.13552
$ws.Range("O5").Value = 0.9210237118384171
$ws.Range("P5").Value = 0.921023711838417
$ws.Range("Q5").Value = 56.26410305678223
$ws.Range("R5").Value = 506.37692751104
$ws.Range("S5").Value = 0.0179302359650433
$ws.Range("T5").Value = 0.0179302359650433

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.160250666666667
$ws.Range("H6").Value = 6.480752000000001
$ws.Range("I6").Value = 0.0194677245922947
$ws.Range("J6").Value = 0.0194677245922947
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3302223333333333
$ws.Range("N6").Value = 0.990667
$ws.Range("O6").Value = 0.01167750336256582
$ws.Range("P6").Value = 0.01167750336256582
$ws.Range("Q6").Value = 0.7133630157315556
$ws.Range("R6").Value = 6.420267141584
$ws.Range("S6").Value = 0.0002273344193880267
$ws.Range("T6").Value = 0.0002273344193880267

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.160250666666667
$ws.Range("H7").Value = 6.480752000000001
$ws.Range("I7").Value = 0.0194677245922947
$ws.Range("J7").Value = 0.0194677245922947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.903109
$ws.Range("N7").Value = 5.709327
$ws.Range("O7").Value = 0.06729878479901708
$ws.Range("P7").Value = 0.06729878479901708
$ws.Range("Q7").Value = 4.111192485989333
$ws.Range("R7").Value = 37.00073237390401
$ws.Range("S7").Value = 0.001310154207863373
$ws.Range("T7").Value = 0.001310154207863373

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd14"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4957383333333333
$ws.Range("H8").Value = 1.487215
$ws.Range("I8").Value = 0.004467489579840358
$ws.Range("J8").Value = 0.004467489579840358
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.04517333333333
$ws.Range("N8").Value = 78.13552
$ws.Range("O8").Value = 0.9210237118384171
$ws.Range("P8").Value = 0.921023711838417
$ws.Range("Q8").Value = 12.91159081964445
$ws.Range("R8").Value = 116.2043173768
$ws.Range("S8").Value = 0.004114663835424017
$ws.Range("T8").Value = 0.004114663835424017

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd14"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4957383333333333
$ws.Range("H9").Value = 1.487215
$ws.Range("I9").Value = 0.004467489579840358
$ws.Range("J9").Value = 0.004467489579840358
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3302223333333333
$ws.Range("N9").Value = 0.990667
$ws.Range("O9").Value = 0.01167750336256582
$ws.Range("P9").Value = 0.01167750336256582
$ws.Range("Q9").Value = 0.1637038691561111
$ws.Range("R9").Value = 1.473334822405
$ws.Range("S9").Value = 0.00005216912459081356
$ws.Range("T9").Value = 0.00005216912459081356

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd14"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4957383333333333
$ws.Range("H10").Value = 1.487215
$ws.Range("I10").Value = 0.004467489579840358
$ws.Range("J10").Value = 0.004467489579840358
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.903109
$ws.Range("N10").Value = 5.709327
$ws.Range("O10").Value = 0.06729878479901708
$ws.Range("P10").Value = 0.06729878479901708
$ws.Range("Q10").Value = 0.9434440838116667
$ws.Range("R10").Value = 8.490996754305
$ws.Range("S10").Value = 0.0003006566198255274
$ws.Range("T10").Value = 0.0003006566198255274
